$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 22; this shifts the old rows 22-27 down to 23-28
$ws.Rows.Item(22).Insert()

# Copy the cell formatting (borders, number format, font) from row 21 into the
# newly inserted row 22 so it matches the other card rows.
$ws.Range("A21:C21").Copy()
$ws.Range("A22:C22").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(22).RowHeight = $ws.Rows.Item(21).RowHeight
$excel.CutCopyMode = 0

# Update the card-type labels and the deposit placeholders for rows 18-22.
$ws.Range("A18").Value = "1、5次2小时卡"
$ws.Range("B18").Value = "{{depositAmount1}}"
$ws.Range("C18").Value = "{{depositAmount1M}}"

$ws.Range("A19").Value = "2、10次全天畅玩卡"
$ws.Range("B19").Value = "{{depositAmount2}}"
$ws.Range("C19").Value = "{{depositAmount2M}}"

$ws.Range("A20").Value = "3、5次亲子2小时卡"
$ws.Range("B20").Value = "{{depositAmount3}}"
$ws.Range("C20").Value = "{{depositAmount3M}}"

$ws.Range("A21").Value = "4、10次亲子全天畅玩卡"
$ws.Range("B21").Value = "{{depositAmount4}}"
$ws.Range("C21").Value = "{{depositAmount4M}}"

$ws.Range("A22").Value = "4、8次情侣畅玩卡"
$ws.Range("B22").Value = "{{depositAmount5}}"
$ws.Range("C22").Value = "{{depositAmount5M}}"

# Keep the view roughly where the author left it.
$ws.Activate()
$ws.Range("C23").Select()
